{"js": "// The document ends with three empty trailing paragraphs. Keep the\n// first one untouched and turn the second/third into the new\n// \"Setup assignment\" notes, matching the target diff.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst count = items.length;\n\n// The last two empty paragraphs at the end of the body.\nconst setupPara = items[count - 2];\nconst filesPara = items[count - 1];\n\nsetupPara.insertText(\"Setup assignment\", \"Replace\");\nfilesPara.insertText(\"Files > project structure\", \"Replace\");\nawait context.sync();\n\n// Two trailing line breaks inside the \"Files > project structure\" paragraph.\nfilesPara.insertBreak(\"Line\", \"End\");\nawait context.sync();\nfilesPara.insertBreak(\"Line\", \"End\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The document ends with three empty trailing paragraphs. Keep the\n# first one untouched and turn the second/third into the new\n# \"Setup assignment\" notes, matching the target diff.\n$count = $d.Paragraphs.Count\n$setupPara = $d.Paragraphs.Item($count - 1)\n$filesPara = $d.Paragraphs.Item($count)\n\n$setupPara.Range.Text = \"Setup assignment\"\n$filesPara.Range.Text = \"Files > project structure\"\n\n# Two trailing line breaks inside the \"Files > project structure\" paragraph.\n$rng = $filesPara.Range\n$rng.Collapse(0)\n$rng.InsertBreak([Microsoft.Office.Interop.Word.WdBreakType]::wdLineBreak)\n$rng.Collapse(0)\n$rng.InsertBreak([Microsoft.Office.Interop.Word.WdBreakType]::wdLineBreak)\n"}
